$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Testcase 1: "Write_Review_Link_Xpath" value text tweak ---
# "Write a review" -> "Write Your Review"
$ws.Range("B11").Value = "Write Your Review"

# Match the formatting used elsewhere for xpath/value cells (blue Courier New),
# copying the format from a cell that already carries it (B29) without
# touching B11's freshly-set text value.
$ws.Range("B29").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# --- Testcase 2: Review_Date_Xpath locator value corrected ---
$ws.Range("B32").Value = "((//li[@class='rvw_title block clear']/div[2]))"

# --- Testcase 7 (in progress): new locator row appended ---
$ws.Range("A32:B32").Copy()
$ws.Range("A33:B33").PasteSpecial(-4122)
$ws.Range("A33").Value = "No_of_Customer_Reviews_DetailsPage_Xpath"
$ws.Range("B33").Value = "//p[@class='rating-links customreview']/a[1]"

# Restore the focused/selected cell as left by the author
$ws.Range("B11").Select()
